$d = $word.ActiveDocument

# --- Step 1: remove the VAR-era _GoBack bookmark that sat after "ALIAS" ---
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks.Item("_GoBack").Delete()
}

# --- Step 2: turn the old "[name]" placeholder into the new "*name" ------
#     shorthand notation (literal match - no wildcards involved).
$find = $d.Content
$find.Find.Execute("[name]", $false, $false, $false, $false, $false, `
                    $false, 1, $false, "*name", 2)

# Re-locate "*name" so we know exactly where the "*" / "name" boundary and
# the trailing end sit, whatever offsets the replace above settled on.
$locate = $d.Content
$locate.Find.Execute("*name", $false, $false, $false, $false, $false, `
                      $false, 1, $false, "", 0)
$starStart = $locate.Start
$starEnd   = $starStart + 1      # boundary between "*" and "name"
$nameEnd   = $starEnd + 4        # boundary right after "name"

# --- Step 3: force clean <w:r> boundaries around "*" and "name" ----------
# Adding then immediately deleting a bookmark at a collapsed point splits
# whatever run currently spans it into two plain runs with no left-over
# formatting, which is exactly how the target markup is structured.
foreach ($pos in @($starStart, $starEnd, $nameEnd)) {
    $splitPoint = $d.Range($pos, $pos)
    $d.Bookmarks.Add("_TmpSplit", $splitPoint)
    $d.Bookmarks.Item("_TmpSplit").Delete()
}

# --- Step 4: plant the _GoBack bookmark right after "name" ---------------
$goBackPoint = $d.Range($nameEnd, $nameEnd)
$d.Bookmarks.Add("_GoBack", $goBackPoint)
